# Apply cryptos list update (Tue Mar 28 07:44:30 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.077.03'
$ws.Range("E2").Value = '  -2.86%  '
$ws.Range("D3").Value = '1.731.80'
$ws.Range("E3").Value = '  -1.39%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.34'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4893'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +6.83%  '
$ws.Range("E8").Value = '  +0.41%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '43.46'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07289'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.94%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.052'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.74%  '
$ws.Range("E12").Value = '  -0.01%  '
$ws.Range("E13").Value = '  -2.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.902'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.23%  '
$ws.Range("D15").Value = '1.731.07'
$ws.Range("E15").Value = '  -1.41%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.905'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.37%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '87.38'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.79%  '
$ws.Range("E18").Value = '  -1.07%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06412'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.63'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.706'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.95%  '
$ws.Range("D23").Value = '27.120.81'
$ws.Range("E23").Value = '  -2.80%  '
$ws.Range("E24").Value = '  -1.78%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.074'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.91%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.99'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.73%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.12%  '
$ws.Range("D28").Value = '1.929.59'
$ws.Range("E28").Value = '  -1.46%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.094'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.32%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '121.50'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.055'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.13%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09361'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.21%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.640'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.78%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.406'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.31%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.06009'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.27%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02193'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.17%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.442'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +7.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '11.03'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.70%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.791'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.14%  '
$ws.Range("E40").Value = '  -2.86%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6023'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.32%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9999'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.101'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.54%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.514'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.36%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.81'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.28%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.584'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5669'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.70%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '119.17'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.40%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.861'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.30%  '
$ws.Range("E50").Value = '  -1.26%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06647'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.24%  '

Write-Host "Applied cryptos list update to" (Get-Date)
